$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276206556/TJaAJ2AJpb-camera-video-segment-driver-1743276209056.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSDJDPAU2S%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T140002Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECIaCXVzLXdlc3QtMiJGMEQCIEEc2n%2FdSkedxvIrawFneB4KmWjRX6zexMgAxdez93GCAiBYlN3X5gu2P5QECbexF8uaZd2CnmIKTdj6eZWt85E0WCrmAwiL%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMKpcBTPh3lRVwmkG9KroDqDa0IZADXHtHrGgsWAtwAVsmtlVoUQIBFw4gHmNHEVmKBsP%2B1KS8OIoHJbBbfWizkpiFhOK5YxeB6NBFbnsQEYdWs0HFvOR%2BmTSvhc%2FmcHozvtxu%2B1z9lcX02DHlYKnq%2B0c3e7zrvqncd0Rnp4YosnlCNqWWkA95AWvIMIsKgKQiaj4mALYtsjyvc71iTQDmv%2BIp0cqOsgPvraW6ZSh5IGb8Rh0PumAtEgpiijT%2BNNuZYHKUaHRJ8%2FeRvjue1UtHoJySjfXZn5aKKLV138wT5gWifoVWgGbVcqk6SMTE5HlIko0dNGM3n3joOi%2BJofybrmMRLeBzpstVQOVfH6jkjRXbmL%2B2raJ%2Btkw%2BI1byq8d8yRLOPZ2qN3UbsAEoRTNun5Se8m%2F1flRT9N5NSmp2kFwyVRafAZljfinfnGFw7ldM8w8NsB4RnR4mHhQTaTlmefw21EIXCzy%2BHEJnmyJQiha2PMw7H0myOOOUZerdhCkSP%2FpRLlsuBkLVdt5F%2Fm4C1T1mOYDqu3ILsl6dn8xT3UQohACBHHNnxjKTPMua4FxcE%2FIDj6oP%2FhrR58tKIX7dz3Dn%2BvBcOBeTujChqqS%2FBjqmAb999puWaBQ2u7GtP4lcrYY4gbGdXkRUDunpLJ%2BgDxYeRuYzuIwMZvhHunaf9QNueXVZ8Ppk%2BGmg4ZK%2Bb5z9MqwwX6y%2FwtAS8NPSsDUNIP193Ti1lZuwVnS8pPUhUtZ6zLO3sfgE6dHxWeHNfdOGzDs7zad8tvr%2FqTW%2BJaF22N8V%2FsGVQypBTbZic9erbj%2FxKaB6Uaur0SERbdxnmNdAjwv472h1li0%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2022%3A00%3A02%20GMT&X-Amz-Signature=67e9ad198569568dddcfe988bcd67507d724696fdd6676e4c6324a0a75cb33c1'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276164863/OF6jgx1Lcp-camera-video-segment-driver-1743276167363.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSDJDPAU2S%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T140002Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECIaCXVzLXdlc3QtMiJGMEQCIEEc2n%2FdSkedxvIrawFneB4KmWjRX6zexMgAxdez93GCAiBYlN3X5gu2P5QECbexF8uaZd2CnmIKTdj6eZWt85E0WCrmAwiL%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMKpcBTPh3lRVwmkG9KroDqDa0IZADXHtHrGgsWAtwAVsmtlVoUQIBFw4gHmNHEVmKBsP%2B1KS8OIoHJbBbfWizkpiFhOK5YxeB6NBFbnsQEYdWs0HFvOR%2BmTSvhc%2FmcHozvtxu%2B1z9lcX02DHlYKnq%2B0c3e7zrvqncd0Rnp4YosnlCNqWWkA95AWvIMIsKgKQiaj4mALYtsjyvc71iTQDmv%2BIp0cqOsgPvraW6ZSh5IGb8Rh0PumAtEgpiijT%2BNNuZYHKUaHRJ8%2FeRvjue1UtHoJySjfXZn5aKKLV138wT5gWifoVWgGbVcqk6SMTE5HlIko0dNGM3n3joOi%2BJofybrmMRLeBzpstVQOVfH6jkjRXbmL%2B2raJ%2Btkw%2BI1byq8d8yRLOPZ2qN3UbsAEoRTNun5Se8m%2F1flRT9N5NSmp2kFwyVRafAZljfinfnGFw7ldM8w8NsB4RnR4mHhQTaTlmefw21EIXCzy%2BHEJnmyJQiha2PMw7H0myOOOUZerdhCkSP%2FpRLlsuBkLVdt5F%2Fm4C1T1mOYDqu3ILsl6dn8xT3UQohACBHHNnxjKTPMua4FxcE%2FIDj6oP%2FhrR58tKIX7dz3Dn%2BvBcOBeTujChqqS%2FBjqmAb999puWaBQ2u7GtP4lcrYY4gbGdXkRUDunpLJ%2BgDxYeRuYzuIwMZvhHunaf9QNueXVZ8Ppk%2BGmg4ZK%2Bb5z9MqwwX6y%2FwtAS8NPSsDUNIP193Ti1lZuwVnS8pPUhUtZ6zLO3sfgE6dHxWeHNfdOGzDs7zad8tvr%2FqTW%2BJaF22N8V%2FsGVQypBTbZic9erbj%2FxKaB6Uaur0SERbdxnmNdAjwv472h1li0%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2022%3A00%3A02%20GMT&X-Amz-Signature=e6a9d43873a80302268de941cdccf572f293ae175404cf6b6190c53e00245042'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276108738/gAxSKLP636-camera-video-segment-driver-1743276111238.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSDJDPAU2S%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T140002Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECIaCXVzLXdlc3QtMiJGMEQCIEEc2n%2FdSkedxvIrawFneB4KmWjRX6zexMgAxdez93GCAiBYlN3X5gu2P5QECbexF8uaZd2CnmIKTdj6eZWt85E0WCrmAwiL%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMKpcBTPh3lRVwmkG9KroDqDa0IZADXHtHrGgsWAtwAVsmtlVoUQIBFw4gHmNHEVmKBsP%2B1KS8OIoHJbBbfWizkpiFhOK5YxeB6NBFbnsQEYdWs0HFvOR%2BmTSvhc%2FmcHozvtxu%2B1z9lcX02DHlYKnq%2B0c3e7zrvqncd0Rnp4YosnlCNqWWkA95AWvIMIsKgKQiaj4mALYtsjyvc71iTQDmv%2BIp0cqOsgPvraW6ZSh5IGb8Rh0PumAtEgpiijT%2BNNuZYHKUaHRJ8%2FeRvjue1UtHoJySjfXZn5aKKLV138wT5gWifoVWgGbVcqk6SMTE5HlIko0dNGM3n3joOi%2BJofybrmMRLeBzpstVQOVfH6jkjRXbmL%2B2raJ%2Btkw%2BI1byq8d8yRLOPZ2qN3UbsAEoRTNun5Se8m%2F1flRT9N5NSmp2kFwyVRafAZljfinfnGFw7ldM8w8NsB4RnR4mHhQTaTlmefw21EIXCzy%2BHEJnmyJQiha2PMw7H0myOOOUZerdhCkSP%2FpRLlsuBkLVdt5F%2Fm4C1T1mOYDqu3ILsl6dn8xT3UQohACBHHNnxjKTPMua4FxcE%2FIDj6oP%2FhrR58tKIX7dz3Dn%2BvBcOBeTujChqqS%2FBjqmAb999puWaBQ2u7GtP4lcrYY4gbGdXkRUDunpLJ%2BgDxYeRuYzuIwMZvhHunaf9QNueXVZ8Ppk%2BGmg4ZK%2Bb5z9MqwwX6y%2FwtAS8NPSsDUNIP193Ti1lZuwVnS8pPUhUtZ6zLO3sfgE6dHxWeHNfdOGzDs7zad8tvr%2FqTW%2BJaF22N8V%2FsGVQypBTbZic9erbj%2FxKaB6Uaur0SERbdxnmNdAjwv472h1li0%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2022%3A00%3A02%20GMT&X-Amz-Signature=f00796e92dbf823cfc49bba3a17e7ad129828556d26762e975dda6ceffc314d8'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205341/1743265652740/hq4516Mzws-camera-video-segment-driver-1743265655240.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSDJDPAU2S%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T140002Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECIaCXVzLXdlc3QtMiJGMEQCIEEc2n%2FdSkedxvIrawFneB4KmWjRX6zexMgAxdez93GCAiBYlN3X5gu2P5QECbexF8uaZd2CnmIKTdj6eZWt85E0WCrmAwiL%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMKpcBTPh3lRVwmkG9KroDqDa0IZADXHtHrGgsWAtwAVsmtlVoUQIBFw4gHmNHEVmKBsP%2B1KS8OIoHJbBbfWizkpiFhOK5YxeB6NBFbnsQEYdWs0HFvOR%2BmTSvhc%2FmcHozvtxu%2B1z9lcX02DHlYKnq%2B0c3e7zrvqncd0Rnp4YosnlCNqWWkA95AWvIMIsKgKQiaj4mALYtsjyvc71iTQDmv%2BIp0cqOsgPvraW6ZSh5IGb8Rh0PumAtEgpiijT%2BNNuZYHKUaHRJ8%2FeRvjue1UtHoJySjfXZn5aKKLV138wT5gWifoVWgGbVcqk6SMTE5HlIko0dNGM3n3joOi%2BJofybrmMRLeBzpstVQOVfH6jkjRXbmL%2B2raJ%2Btkw%2BI1byq8d8yRLOPZ2qN3UbsAEoRTNun5Se8m%2F1flRT9N5NSmp2kFwyVRafAZljfinfnGFw7ldM8w8NsB4RnR4mHhQTaTlmefw21EIXCzy%2BHEJnmyJQiha2PMw7H0myOOOUZerdhCkSP%2FpRLlsuBkLVdt5F%2Fm4C1T1mOYDqu3ILsl6dn8xT3UQohACBHHNnxjKTPMua4FxcE%2FIDj6oP%2FhrR58tKIX7dz3Dn%2BvBcOBeTujChqqS%2FBjqmAb999puWaBQ2u7GtP4lcrYY4gbGdXkRUDunpLJ%2BgDxYeRuYzuIwMZvhHunaf9QNueXVZ8Ppk%2BGmg4ZK%2Bb5z9MqwwX6y%2FwtAS8NPSsDUNIP193Ti1lZuwVnS8pPUhUtZ6zLO3sfgE6dHxWeHNfdOGzDs7zad8tvr%2FqTW%2BJaF22N8V%2FsGVQypBTbZic9erbj%2FxKaB6Uaur0SERbdxnmNdAjwv472h1li0%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2022%3A00%3A02%20GMT&X-Amz-Signature=45e0f29877ab2bab84664f3660a05f052d234fb593ca440331250add52e71b9d'
$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743255999969/Yjq7ZuVo4p-camera-video-segment-driver-1743256002469.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSDJDPAU2S%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T140002Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECIaCXVzLXdlc3QtMiJGMEQCIEEc2n%2FdSkedxvIrawFneB4KmWjRX6zexMgAxdez93GCAiBYlN3X5gu2P5QECbexF8uaZd2CnmIKTdj6eZWt85E0WCrmAwiL%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMKpcBTPh3lRVwmkG9KroDqDa0IZADXHtHrGgsWAtwAVsmtlVoUQIBFw4gHmNHEVmKBsP%2B1KS8OIoHJbBbfWizkpiFhOK5YxeB6NBFbnsQEYdWs0HFvOR%2BmTSvhc%2FmcHozvtxu%2B1z9lcX02DHlYKnq%2B0c3e7zrvqncd0Rnp4YosnlCNqWWkA95AWvIMIsKgKQiaj4mALYtsjyvc71iTQDmv%2BIp0cqOsgPvraW6ZSh5IGb8Rh0PumAtEgpiijT%2BNNuZYHKUaHRJ8%2FeRvjue1UtHoJySjfXZn5aKKLV138wT5gWifoVWgGbVcqk6SMTE5HlIko0dNGM3n3joOi%2BJofybrmMRLeBzpstVQOVfH6jkjRXbmL%2B2raJ%2Btkw%2BI1byq8d8yRLOPZ2qN3UbsAEoRTNun5Se8m%2F1flRT9N5NSmp2kFwyVRafAZljfinfnGFw7ldM8w8NsB4RnR4mHhQTaTlmefw21EIXCzy%2BHEJnmyJQiha2PMw7H0myOOOUZerdhCkSP%2FpRLlsuBkLVdt5F%2Fm4C1T1mOYDqu3ILsl6dn8xT3UQohACBHHNnxjKTPMua4FxcE%2FIDj6oP%2FhrR58tKIX7dz3Dn%2BvBcOBeTujChqqS%2FBjqmAb999puWaBQ2u7GtP4lcrYY4gbGdXkRUDunpLJ%2BgDxYeRuYzuIwMZvhHunaf9QNueXVZ8Ppk%2BGmg4ZK%2Bb5z9MqwwX6y%2FwtAS8NPSsDUNIP193Ti1lZuwVnS8pPUhUtZ6zLO3sfgE6dHxWeHNfdOGzDs7zad8tvr%2FqTW%2BJaF22N8V%2FsGVQypBTbZic9erbj%2FxKaB6Uaur0SERbdxnmNdAjwv472h1li0%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2022%3A00%3A02%20GMT&X-Amz-Signature=f19a2600aae6e2c2166aef0eff1dc8e5542b1fc14665c185b2f32d8fe0fa2d61'
